$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 30
$ws.Range("I6").Value = 30
$ws.Range("K6").Value = 90
$ws.Range("M6").Value = 22
$ws.Range("H9").Value = 209.53334
$ws.Range("I9").Value = 147.48
$ws.Range("J9").Value = 519.8
$ws.Range("K9").Value = 147.48
$ws.Range("L9").Value = 519.8
$ws.Range("M9").Value = 21.52000000000001
$ws.Range("N9").Value = -857.8
$ws.Range("H28").Value = 1470.4348
$ws.Range("I28").Value = 1095.6875
$ws.Range("J28").Value = 2327
$ws.Range("K28").Value = 1095.6875
$ws.Range("L28").Value = 2327
$ws.Range("M28").Value = -610.6875
$ws.Range("N28").Value = -3297
$ws.Range("H40").Value = 6254006
$ws.Range("I40").Value = 16669500
$ws.Range("J40").Value = 4709.9
$ws.Range("K40").Value = 16669500
$ws.Range("L40").Value = 4709.9
$ws.Range("M40").Value = -16669325
$ws.Range("N40").Value = -5059.9
$ws.Range("H41").Value = 533.913
$ws.Range("I41").Value = 394.91666
$ws.Range("J41").Value = 685.5454999999999
$ws.Range("K41").Value = 394.91666
$ws.Range("L41").Value = 685.5454999999999
$ws.Range("M41").Value = 45.08334000000002
$ws.Range("N41").Value = -1565.5455
$ws.Range("H62").Value = 123980.31
$ws.Range("I62").Value = 140834.64
$ws.Range("K62").Value = 140834.64
$ws.Range("M62").Value = -140210.64
$ws.Range("H65").Value = 123980.31
$ws.Range("I65").Value = 140834.64
$ws.Range("K65").Value = 704173.2000000001
$ws.Range("M65").Value = -701053.2000000001
$ws.Range("H86").Value = 950.3333
$ws.Range("J86").Value = 1000
$ws.Range("L86").Value = 1000
$ws.Range("N86").Value = -3246
$ws.Range("H87").Value = 70354
$ws.Range("J87").Value = 70354
$ws.Range("L87").Value = 70354
$ws.Range("N87").Value = -72850
$ws.Range("H89").Value = 950.3333
$ws.Range("J89").Value = 1000
$ws.Range("L89").Value = 5000
$ws.Range("N89").Value = -16232
$ws.Range("H90").Value = 70354
$ws.Range("J90").Value = 70354
$ws.Range("L90").Value = 211062
$ws.Range("N90").Value = -223542
$ws.Range("H92").Value = 243.8
$ws.Range("I92").Value = 303.33334
$ws.Range("J92").Value = 5.6666665
$ws.Range("K92").Value = 303.33334
$ws.Range("L92").Value = 5.6666665
$ws.Range("M92").Value = 944.66666
$ws.Range("N92").Value = -2501.6666665
$ws.Range("H112").Value = 837282.4399999999
$ws.Range("I112").Value = 3050
$ws.Range("J112").Value = 1254398.6
$ws.Range("K112").Value = 9150
$ws.Range("L112").Value = 3763195.8
$ws.Range("M112").Value = -8042
$ws.Range("N112").Value = -3765411.8
$ws.Range("H116").Value = 3672.1428
$ws.Range("I116").Value = 3650
$ws.Range("J116").Value = 3727.5
$ws.Range("K116").Value = 3650
$ws.Range("L116").Value = 3727.5
$ws.Range("M116").Value = -208
$ws.Range("N116").Value = -10611.5
$ws.Range("H132").Value = 37043572
$ws.Range("I132").Value = 47625304
$ws.Range("K132").Value = 142875912
$ws.Range("M132").Value = -142873382
$ws.Range("H135").Value = 6172.839
$ws.Range("I135").Value = 4262.75
$ws.Range("J135").Value = 12721.714
$ws.Range("K135").Value = 38364.75
$ws.Range("L135").Value = 114495.426
$ws.Range("M135").Value = -35829.75
$ws.Range("N135").Value = -119565.426
$ws.Range("H137").Value = 3719.2666
$ws.Range("J137").Value = 3413.5715
$ws.Range("L137").Value = 10240.7145
$ws.Range("N137").Value = -15340.7145
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 950
$ws.Range("I16").Value = 900
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 900
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -613
$ws.Range("N16").Value = -1574
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H21").Value = 15395
$ws.Range("I21").Value = 15395
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 15395
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -15021
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H25").Value = 1858
$ws.Range("I25").Value = 1858
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1858
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -1456
$ws.Range("N25").ClearContents()
$ws.Range("H28").Value = 17700
$ws.Range("J28").Value = 52759.5
$ws.Range("L28").Value = 52759.5
$ws.Range("N28").Value = -53143.5
$ws.Range("H30").Value = 1900
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H32").Value = 2147.8372
$ws.Range("I32").Value = 2098.2354
$ws.Range("K32").Value = 2098.2354
$ws.Range("M32").Value = -1811.2354
$ws.Range("H61").Value = 4698.3784
$ws.Range("I61").Value = 5047.7666
$ws.Range("K61").Value = 5047.7666
$ws.Range("M61").Value = -4835.7666
$ws.Range("H63").Value = 2944.6667
$ws.Range("I63").Value = 2725.6667
$ws.Range("K63").Value = 2725.6667
$ws.Range("M63").Value = -2039.6667
$ws.Range("H66").Value = 2944.6667
$ws.Range("I66").Value = 2725.6667
$ws.Range("K66").Value = 13628.3335
$ws.Range("M66").Value = -10196.3335
$ws.Range("H99").Value = 17700
$ws.Range("J99").Value = 52759.5
$ws.Range("L99").Value = 52759.5
$ws.Range("N99").Value = -58749.5
$ws.Range("H102").Value = 40000710
$ws.Range("I102").Value = 1279.5
$ws.Range("J102").Value = 66667000
$ws.Range("K102").Value = 1279.5
$ws.Range("L102").Value = 66667000
$ws.Range("M102").Value = 342.5
$ws.Range("N102").Value = -66670244
$ws.Range("H122").Value = 2387.4119
$ws.Range("I122").Value = 2224.125
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6672.375
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -4222.375
$ws.Range("N122").Value = -19900
$ws.Range("H136").Value = 4698.3784
$ws.Range("I136").Value = 5047.7666
$ws.Range("K136").Value = 15143.2998
$ws.Range("M136").Value = -12593.2998
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 63332.668
$ws.Range("J2").Value = 64999
$ws.Range("L2").Value = 64999
$ws.Range("N2").Value = -65225
$ws.Range("H6").Value = 60000
$ws.Range("I6").Value = 50000
$ws.Range("K6").Value = 50000
$ws.Range("M6").Value = -49887
$ws.Range("H11").Value = 2650.3076
$ws.Range("I11").Value = 2839.4443
$ws.Range("J11").Value = 2224.75
$ws.Range("K11").Value = 2839.4443
$ws.Range("L11").Value = 2224.75
$ws.Range("M11").Value = -2699.4443
$ws.Range("N11").Value = -2504.75
$ws.Range("H20").Value = 1027.7391
$ws.Range("I20").Value = 1045.3077
$ws.Range("J20").Value = 1004.9
$ws.Range("K20").Value = 1045.3077
$ws.Range("L20").Value = 1004.9
$ws.Range("M20").Value = -798.3077000000001
$ws.Range("N20").Value = -1498.9
$ws.Range("H28").Value = 64330.668
$ws.Range("J28").Value = 66250
$ws.Range("L28").Value = 66250
$ws.Range("N28").Value = -66838
$ws.Range("H80").Value = 474.1905
$ws.Range("J80").Value = 500.35715
$ws.Range("L80").Value = 500.35715
$ws.Range("N80").Value = -2496.35715
$ws.Range("H83").Value = 474.1905
$ws.Range("J83").Value = 500.35715
$ws.Range("L83").Value = 2501.78575
$ws.Range("N83").Value = -12485.78575
$ws.Range("H86").Value = 2909859.8
$ws.Range("I86").Value = 4458143
$ws.Range("J86").Value = 6828.625
$ws.Range("K86").Value = 4458143
$ws.Range("L86").Value = 6828.625
$ws.Range("M86").Value = -4457020
$ws.Range("N86").Value = -9074.625
$ws.Range("H89").Value = 2909859.8
$ws.Range("I89").Value = 4458143
$ws.Range("J89").Value = 6828.625
$ws.Range("K89").Value = 22290715
$ws.Range("L89").Value = 34143.125
$ws.Range("M89").Value = -22285099
$ws.Range("N89").Value = -45375.125
$ws.Range("H99").Value = 2460.1177
$ws.Range("I99").Value = 1726.1666
$ws.Range("K99").Value = 1726.1666
$ws.Range("M99").Value = -228.1666
$ws.Range("H105").Value = 2832.22
$ws.Range("I105").Value = 1563.6666
$ws.Range("K105").Value = 1563.6666
$ws.Range("M105").Value = 183.3334
$ws.Range("H118").Value = 65000
$ws.Range("J118").Value = 65000
$ws.Range("L118").Value = 65000
$ws.Range("N118").Value = -68314
$ws.Range("H134").Value = 3569.48
$ws.Range("I134").Value = 3655.7083
$ws.Range("K134").Value = 10967.1249
$ws.Range("M134").Value = -8432.124899999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1126.4783
$ws.Range("J58").Value = 1539
$ws.Range("L58").Value = 1539
$ws.Range("N58").Value = -1945
$ws.Range("H107").Value = 407.94116
$ws.Range("I107").Value = 275.83334
$ws.Range("K107").Value = 275.83334
$ws.Range("M107").Value = 1644.16666
$ws.Range("H122").Value = 3370.5386
$ws.Range("I122").Value = 2031.9
$ws.Range("K122").Value = 6095.700000000001
$ws.Range("M122").Value = -3645.700000000001
$ws.Range("H134").Value = 1595.5625
$ws.Range("J134").Value = 2022.25
$ws.Range("L134").Value = 6066.75
$ws.Range("N134").Value = -11136.75
$ws.Range("H136").Value = 1126.4783
$ws.Range("J136").Value = 1539
$ws.Range("L136").Value = 4617
$ws.Range("N136").Value = -9717
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 169816
$ws.Range("I8").Value = 169816
$ws.Range("K8").Value = 509448
$ws.Range("M8").Value = -509309
$ws.Range("H11").Value = 34.666668
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 34.666668
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 104.000004
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -384.000004
$ws.Range("H34").Value = 2138.6667
$ws.Range("I34").Value = 1735.7646
$ws.Range("J34").Value = 8988
$ws.Range("K34").Value = 5207.293799999999
$ws.Range("L34").Value = 26964
$ws.Range("M34").Value = -5123.293799999999
$ws.Range("N34").Value = -27132
$ws.Range("H38").Value = 202.26315
$ws.Range("J38").Value = 174.5
$ws.Range("L38").Value = 523.5
$ws.Range("N38").Value = -1217.5
$ws.Range("H113").Value = 956.3
$ws.Range("I113").Value = 661.6667
$ws.Range("K113").Value = 1985.0001
$ws.Range("M113").Value = 184.9999
$ws.Range("H114").Value = 5500
$ws.Range("I114").Value = 1000
$ws.Range("J114").Value = 10000
$ws.Range("K114").Value = 3000
$ws.Range("L114").Value = 30000
$ws.Range("M114").Value = 254
$ws.Range("N114").Value = -36508
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H117").Value = 1768.3334
$ws.Range("J117").Value = 2022
$ws.Range("L117").Value = 6066
$ws.Range("N117").Value = -12950
$ws.Range("H120").Value = 9275.666999999999
$ws.Range("I120").Value = 9275.666999999999
$ws.Range("K120").Value = 27827.001
$ws.Range("M120").Value = -22989.001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 16600
$ws.Range("H92").Value = 20551.75
$ws.Range("J92").Value = 20551.75
$ws.Range("L92").Value = 20551.75
$ws.Range("N92").Value = -24295.75
$ws.Range("H99").Value = 20700.363
$ws.Range("I99").Value = 2022.6
$ws.Range("K99").Value = 2022.6
$ws.Range("M99").Value = 223.4000000000001
$ws.Range("H123").Value = 50305.637
$ws.Range("J123").Value = 50305.637
$ws.Range("L123").Value = 50305.637
$ws.Range("N123").Value = -55205.637
$ws.Range("H126").Value = 32438.53
$ws.Range("I126").Value = 3782.125
$ws.Range("J126").Value = 57910.89
$ws.Range("K126").Value = 11346.375
$ws.Range("L126").Value = 173732.67
$ws.Range("M126").Value = -8876.375
$ws.Range("N126").Value = -178672.67
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 6936.2856
$ws.Range("J4").Value = 5000
$ws.Range("L4").Value = 5000
$ws.Range("N4").Value = -5226
$ws.Range("H7").Value = 3725.1428
$ws.Range("I7").Value = 2916.9443
$ws.Range("K7").Value = 2916.9443
$ws.Range("M7").Value = -2804.9443
$ws.Range("H22").Value = 6815.154
$ws.Range("I22").Value = 2759.4
$ws.Range("K22").Value = 2759.4
$ws.Range("M22").Value = -2464.4
$ws.Range("H27").Value = 6815.154
$ws.Range("I27").Value = 2759.4
$ws.Range("K27").Value = 2759.4
$ws.Range("M27").Value = -2652.4
$ws.Range("H28").Value = 6936.2856
$ws.Range("J28").Value = 5000
$ws.Range("L28").Value = 5000
$ws.Range("N28").Value = -5464
$ws.Range("H29").Value = 14016
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H37").Value = 6936.2856
$ws.Range("J37").Value = 5000
$ws.Range("L37").Value = 5000
$ws.Range("N37").Value = -5214
$ws.Range("H46").Value = 14010.546
$ws.Range("I46").Value = 740
$ws.Range("J46").Value = 25069.334
$ws.Range("K46").Value = 740
$ws.Range("L46").Value = 25069.334
$ws.Range("M46").Value = -552
$ws.Range("N46").Value = -25445.334
$ws.Range("H93").Value = 1420.5714
$ws.Range("I93").Value = 1050.75
$ws.Range("J93").Value = 1568.5
$ws.Range("K93").Value = 1050.75
$ws.Range("L93").Value = 1568.5
$ws.Range("M93").Value = 197.25
$ws.Range("N93").Value = -4064.5
$ws.Range("H122").Value = 5029.684
$ws.Range("I122").Value = 4755.357
$ws.Range("J122").Value = 5797.8
$ws.Range("K122").Value = 14266.071
$ws.Range("L122").Value = 17393.4
$ws.Range("M122").Value = -11816.071
$ws.Range("N122").Value = -22293.4
$ws.Range("H126").Value = 3725.1428
$ws.Range("I126").Value = 2916.9443
$ws.Range("K126").Value = 8750.832900000001
$ws.Range("M126").Value = -6280.832900000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 17509.5
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 17509.5
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H100").Value = 1329.8889
$ws.Range("I100").Value = 1264.6428
$ws.Range("J100").Value = 1558.25
$ws.Range("K100").Value = 2529.2856
$ws.Range("L100").Value = 3116.5
$ws.Range("M100").Value = -1988.2856
$ws.Range("N100").Value = -4198.5
$ws.Range("H107").Value = 639.44446
$ws.Range("I107").Value = 560.0714
$ws.Range("J107").Value = 724.9231
$ws.Range("K107").Value = 1680.2142
$ws.Range("L107").Value = 2174.7693
$ws.Range("M107").Value = 239.7857999999999
$ws.Range("N107").Value = -6014.7693
$ws.Range("H122").Value = 1922.7826
$ws.Range("I122").Value = 2065.1707
$ws.Range("K122").Value = 6195.5121
$ws.Range("M122").Value = -3745.5121
$ws.Range("H126").Value = 3993.1072
$ws.Range("I126").Value = 3993.1072
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 11979.3216
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9509.321599999999
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2756.5715
$ws.Range("I132").Value = 2859.4
$ws.Range("J132").Value = 2499.5
$ws.Range("K132").Value = 8578.200000000001
$ws.Range("L132").Value = 7498.5
$ws.Range("M132").Value = -6048.200000000001
$ws.Range("N132").Value = -12558.5
$ws.Range("H136").Value = 3288.3547
$ws.Range("I136").Value = 2711.8572
$ws.Range("J136").Value = 4499
$ws.Range("K136").Value = 8135.571599999999
$ws.Range("L136").Value = 13497
$ws.Range("M136").Value = -5585.571599999999
$ws.Range("N136").Value = -18597
